$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (251) down to the new rows (252:255)
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)

$data = @(
    @(252, 44326, 10, 66, 91.19927040583676),
    @(253, 44327, 10, 72, 99.49011317000372),
    @(254, 44328, 0, 71, 98.10830604264255),
    @(255, 44329, 14, 76, 105.0173416794484)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
